$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 1: in "496、503、456、316、402、321、84、85" paragraph,
# split the trailing run "、456、316、402、321、84、85" into three
# runs, coloring "、456" and "、316、" green (00B050) while leaving
# "402、321、84、85" with no explicit color.
# ---------------------------------------------------------------
$green = 5287936  # RGB(0x00,0xB0,0x50) == 00B050 as a Word BGR-packed long

$r1 = $d.Content
$r1.Find.Execute("、456", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.Font.Color = $green

$r2 = $d.Content
$r2.Find.Execute("、316、", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Font.Color = $green

# ---------------------------------------------------------------
# Change 2: in "89、136、137、260、268" paragraph, split the run
# into three runs with identical (default) formatting:
#   "89、136、137、"  "260"  "、268"
# A temporary bookmark is added/removed around each boundary
# substring; this forces the run to split at that exact position
# without touching any character formatting (no stray <w:rPr/>).
# ---------------------------------------------------------------
$b1 = $d.Content
$b1.Find.Execute("260", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("zzTmpSplit1", $b1)
$d.Bookmarks("zzTmpSplit1").Delete()

$b2 = $d.Content
$b2.Find.Execute("、268", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("zzTmpSplit2", $b2)
$d.Bookmarks("zzTmpSplit2").Delete()
